# Update the BDD "First Name" step text (drop "Please ensure " from the sentence)
# per the new feature-step format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 'When you correctly enter "First Name" before moving to the next field.'

# Remove the hyperlink that was attached to F2 (the form URL cell).
[void]$ws.Hyperlinks.Delete()

# Move/restore the active selection to C3.
[void]$ws.Range("C3").Select()
